$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 17298
$ws.Range("C3").Value = 2691
$ws.Range("D3").Value = 2982
$ws.Range("B4").Value = 8543
$ws.Range("C4").Value = 1106
$ws.Range("D4").Value = 1127
$ws.Range("B5").Value = 30307
$ws.Range("C5").Value = 2643
$ws.Range("D5").Value = 2941
$ws.Range("B6").Value = 425
$ws.Range("C6").Value = 283
$ws.Range("D6").Value = 67
$ws.Range("B7").Value = 18753
$ws.Range("C7").Value = 3426
$ws.Range("D7").Value = 2983
$ws.Range("B8").Value = 2345
$ws.Range("C8").Value = 603
$ws.Range("D8").Value = 571
$ws.Range("B9").Value = 2396
$ws.Range("C9").Value = 462
$ws.Range("D9").Value = 307
$ws.Range("B10").Value = 935
$ws.Range("C10").Value = 185
$ws.Range("D10").Value = 120
$ws.Range("B11").Value = 61
$ws.Range("C11").Value = 103
$ws.Range("D11").Value = 17
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("B13").Value = 438
$ws.Range("C13").Value = 114
$ws.Range("D13").Value = 149
$ws.Range("B14").Value = 1260
$ws.Range("C14").Value = 583
$ws.Range("D14").Value = 372
$ws.Range("B15").Value = 2195
$ws.Range("C15").Value = 835
$ws.Range("D15").Value = 362
$ws.Range("B16").Value = 1246
$ws.Range("C16").Value = 538
$ws.Range("D16").Value = 170
$ws.Range("B17").Value = 612
$ws.Range("C17").Value = 307
$ws.Range("D17").Value = 58
$ws.Range("B18").Value = 6514
$ws.Range("C18").Value = 1265
$ws.Range("D18").Value = 1375
$ws.Range("B19").Value = 731
$ws.Range("C19").Value = 300
$ws.Range("D19").Value = 269
$ws.Range("B20").Value = 7888
$ws.Range("C20").Value = 911
$ws.Range("D20").Value = 1359
$ws.Range("B21").Value = 102
$ws.Range("C21").Value = 164
$ws.Range("D21").Value = 7
$ws.Range("B22").Value = 6924
$ws.Range("C22").Value = 974
$ws.Range("D22").Value = 1400
$ws.Range("B23").Value = 566
$ws.Range("C23").Value = 208
$ws.Range("D23").Value = 58
$ws.Range("B24").Value = 7517
$ws.Range("C24").Value = 1284
$ws.Range("D24").Value = 1635
$ws.Range("B25").Value = 33005
$ws.Range("C25").Value = 3228
$ws.Range("D25").Value = 4152
$ws.Range("B26").Value = 2250
$ws.Range("C26").Value = 790
$ws.Range("D26").Value = 444
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("B28").Value = 2241
$ws.Range("C28").Value = 438
$ws.Range("D28").Value = 592
$ws.Range("B29").Value = 527
$ws.Range("C29").Value = 216
$ws.Range("D29").Value = 115
$ws.Range("B30").Value = 6015
$ws.Range("C30").Value = 1171
$ws.Range("D30").Value = 1039
$ws.Range("B31").Value = 236
$ws.Range("C31").Value = 77
$ws.Range("D31").Value = 132
$ws.Range("B32").Value = 970
$ws.Range("C32").Value = 750
$ws.Range("D32").Value = 191
$ws.Range("B33").Value = 6444
$ws.Range("C33").Value = 1581
$ws.Range("D33").Value = 1183
$ws.Range("B34").Value = 4044
$ws.Range("C34").Value = 1316
$ws.Range("D34").Value = 1013
$ws.Range("B35").Value = 2547
$ws.Range("C35").Value = 308
$ws.Range("D35").Value = 691
$ws.Range("B36").Value = 22245
$ws.Range("C36").Value = 2495
$ws.Range("D36").Value = 2378
$ws.Range("B37").Value = 3148
$ws.Range("C37").Value = 1328
$ws.Range("D37").Value = 546
$ws.Range("B38").Value = 9657
$ws.Range("C38").Value = 908
$ws.Range("D38").Value = 1240
$ws.Range("B39").Value = 394
$ws.Range("C39").Value = 413
$ws.Range("D39").Value = 142
$ws.Range("B40").Value = 885
$ws.Range("C40").Value = 236
$ws.Range("D40").Value = 301
$ws.Range("B41").Value = 1661
$ws.Range("C41").Value = 203
$ws.Range("D41").Value = 88
$ws.Range("B42").Value = 6333
$ws.Range("C42").Value = 361
$ws.Range("D42").Value = 189
$ws.Range("B43").Value = 190
$ws.Range("C43").Value = 53
$ws.Range("D43").Value = 62
$ws.Range("B44").Value = 411
$ws.Range("C44").Value = 34
$ws.Range("D44").Value = 38
$ws.Range("B45").Value = 851
$ws.Range("C45").Value = 14
$ws.Range("D45").Value = 1
$ws.Range("B46").Value = 1180
$ws.Range("C46").Value = 401
$ws.Range("D46").Value = 169
$ws.Range("B47").Value = 4599
$ws.Range("C47").Value = 1529
$ws.Range("D47").Value = 906
$ws.Range("B48").Value = 12499
$ws.Range("C48").Value = 1681
$ws.Range("D48").Value = 1997
$ws.Range("B49").Value = 5210
$ws.Range("C49").Value = 1524
$ws.Range("D49").Value = 545
$ws.Range("B50").Value = 4749
$ws.Range("C50").Value = 467
$ws.Range("D50").Value = 666
$ws.Range("B51").Value = 12254
$ws.Range("C51").Value = 1283
$ws.Range("D51").Value = 1926
$ws.Range("B52").Value = 1859
$ws.Range("C52").Value = 249
$ws.Range("D52").Value = 537
$ws.Range("B53").Value = 6154
$ws.Range("C53").Value = 1350
$ws.Range("D53").Value = 1134
$ws.Range("B54").Value = 646
$ws.Range("C54").Value = 485
$ws.Range("D54").Value = 275
$ws.Range("B55").Value = 890
$ws.Range("C55").Value = 606
$ws.Range("D55").Value = 75
$ws.Range("B56").Value = 975
$ws.Range("C56").Value = 256
$ws.Range("D56").Value = 299
$ws.Range("B57").Value = 5139
$ws.Range("C57").Value = 2070
$ws.Range("D57").Value = 1053
$ws.Range("B58").Value = 9446
$ws.Range("C58").Value = 665
$ws.Range("D58").Value = 378
$ws.Range("B59").Value = 272128
$ws.Range("C59").Value = 44736
$ws.Range("D59").Value = 41395
